$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: mark Wednesday/Column E as worked (copy style from D4, then set value) ---
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 1

# --- Shift the summary block in D28:E33 one column to the right (D->E, E->F) ---
# Do the right-most move first so we do not overwrite data we still need.
$ws.Range("E28:E33").Cut($ws.Range("F28:F33"))
$ws.Range("D28:D33").Cut($ws.Range("E28:E33"))

# --- Add the new "Uren" column (C) ---
$ws.Range("C28").Value = "Uren"
$ws.Range("C30").Value = 8.5

# --- Restore the formulas (Cut turns formulas into static values), with updated references ---
$ws.Range("F29").Formula = "=SUM(C2:G26)"
$ws.Range("F30").Formula = "=F29*C30"
$ws.Range("F31").Formula = "=F28-F30"
$ws.Range("F32").Formula = "=F31/C30"
$ws.Range("F33").Formula = "=F32/4"

# --- Column widths: D shrinks, E becomes a custom width ---
$ws.Columns("D").ColumnWidth = 12.166666666666666
$ws.Columns("E").ColumnWidth = 15.666666666666666

# --- Update selection to match the last edited cell ---
$ws.Range("F33").Select()
